$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new data points (rows 10 and 11) to the Standard/Ratio table.
$ws.Range("A10").Value = 15
$ws.Range("B10").Value = 14.676
$ws.Range("A11").Value = 20
$ws.Range("B11").Value = 19.4

# Match the formatting used for the rest of column A: centered, with a thin
# left/right border (column B keeps inheriting its existing column style).
$a10 = $ws.Range("A10")
$a10.Borders.Item(7).LineStyle = 1
$a10.Borders.Item(10).LineStyle = 1
$a10.HorizontalAlignment = -4108

$a10.Copy()
$ws.Range("A11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Leave the selection where the user would land after typing the last value.
$ws.Range("B12").Select()
